# edit.ps1 - apply sprint_tasks_codex.xlsx S21 update
#
# Notes on this runtime's PowerShell-COM quirks (discovered empirically):
#  - Named parameters on user-defined functions (e.g. "-foo bar") do not bind;
#    use strictly positional parameters for custom functions.
#  - Range "get value" must be called as a method: $range.Value() not $range.Value
#  - Range "set value" with a 2D payload needs a real .NET 2D object[,] array
#    (jagged PowerShell @() literals do not marshal correctly); and the target
#    Range must be obtained via a single "A1:H1"-style address string.

function Set-RowValues($ws, $rowNum, $values) {
    $n = $values.Count
    $arr2d = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr2d[0, $i] = $values[$i]
    }
    $lastColLetter = [char](64 + $n)
    $addr = "A" + $rowNum + ":" + $lastColLetter + $rowNum
    $rng = $ws.Range($addr)
    $rng.Value = $arr2d
}

function Copy-RowFormat($ws, $srcRow, $dstRow) {
    $srcAddr = "A" + $srcRow + ":H" + $srcRow
    $dstAddr = "A" + $dstRow + ":H" + $dstRow
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Application.CutCopyMode = $false
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 167-173: row height tweak 41.75 -> 41.25 (content/style unchanged)
# ---------------------------------------------------------------------------
for ($r = 167; $r -le 173; $r++) {
    $ws.Rows.Item($r).RowHeight = 41.25
}

# ---------------------------------------------------------------------------
# 2) Row 174: style it like the rest of the table (it was previously plain)
#    and update status/remarks to reflect implementation.
# ---------------------------------------------------------------------------
Copy-RowFormat $ws 167 174
$ws.Rows.Item(174).RowHeight = 41.75
Set-RowValues $ws 174 @(
    "S21",
    "G01",
    "Custom bracket orders (primary + follow-up GTT legs)",
    "S21_G01_TB001",
    "Design and implement backend helper/flow to create paired manual orders (primary + LIMIT GTT leg) given side, qty, effective price, and MTP.",
    "Keeps invariant that each leg is a normal WAITING manual order; does not introduce broker-specific bracket types.",
    "implemented",
    "Bracket pairing implemented by creating a second LIMIT+GTT manual order from the frontend using existing /api/orders/, avoiding extra backend helper complexity in this phase."
)

# ---------------------------------------------------------------------------
# 3) Row 175: new task S21_G01_TF001 (Bracket section in Holdings dialog)
# ---------------------------------------------------------------------------
Copy-RowFormat $ws 167 175
$ws.Rows.Item(175).RowHeight = 55.2
Set-RowValues $ws 175 @(
    "S21",
    "G01",
    "Custom bracket orders (primary + follow-up GTT legs)",
    "S21_G01_TF001",
    "Extend Holdings Buy/Sell dialog with a Bracket section that lets the user enable a follow-up GTT leg, pre-fills MTP% from current appreciation, and previews the derived target price.",
    "Bracket invocation remains manual; dialog simply creates two manual orders via the existing orders API.",
    "implemented",
    "Holdings Buy/Sell dialog now has a Bracket section with enable checkbox, MTP% field, and live GTT price preview."
)

# ---------------------------------------------------------------------------
# 4) Row 176: new task S21_G01_TF002 (Queue/Orders grid bracket highlighting)
# ---------------------------------------------------------------------------
Copy-RowFormat $ws 167 176
$ws.Rows.Item(176).RowHeight = 41.75
Set-RowValues $ws 176 @(
    "S21",
    "G01",
    "Custom bracket orders (primary + follow-up GTT legs)",
    "S21_G01_TF002",
    "Highlight bracket-related information in the Queue and Orders grids (order_type, trigger_price, GTT flag) and optionally tag bracket legs for easier identification.",
    "Builds on the existing DataGrid-based queue and orders views added earlier.",
    "implemented",
    "Queue and Orders grids converted to DataGrid and now expose order_type, trigger_price, and GTT flags clearly (e.g. LIMIT (GTT))."
)

# ---------------------------------------------------------------------------
# 5) Row 177: task S21_G02_TB001 moves down from its old row 177 (still
#    "planned", unchanged content, but picks up table row styling + new ht).
# ---------------------------------------------------------------------------
Copy-RowFormat $ws 167 177
$ws.Rows.Item(177).RowHeight = 55.2
Set-RowValues $ws 177 @(
    "S21",
    "G02",
    "Bracket-order backtesting using Kite OHLCV",
    "S21_G02_TB001",
    "Add a backend console script that pulls OHLCV via the existing market-data layer and simulates the custom bracket logic over a given symbol, timeframe, and lookback.",
    "Focus initial experiments on BSE and NETWEB over the last month but keep the script parameterised for any symbol.",
    "planned",
    "Enables quantitative evaluation of how often the MTP-based bracket legs would have been filled and the resulting P&L and drawdowns."
)

# ---------------------------------------------------------------------------
# 6) Row 178: brand new task S21_G01_TF003 (SELL-side default MTP refinement)
# ---------------------------------------------------------------------------
Copy-RowFormat $ws 167 178
$ws.Rows.Item(178).RowHeight = 41.75
Set-RowValues $ws 178 @(
    "S21",
    "G01",
    "Custom bracket orders (primary + follow-up GTT legs)",
    "S21_G01_TF003",
    "Refine SELL-side default MTP logic so that it mirrors positive Today P&L% only above a small threshold and clamps to a reasonable min/max band.",
    "Implements a 3% appreciation threshold with 3–20% clamps instead of always using raw gain.",
    "implemented",
    "Uses today_pnl_percent from holdings with a 3% threshold and 3–20% clamps, instead of raw gain vs average price."
)

# ---------------------------------------------------------------------------
# 7) Row 179: brand new task S21_G03_TD001, kept plain (no special row
#    height / fill style), matching the rest of the appended rows pattern.
#    A brand-new row implicitly inherits the wrap-text column formatting as
#    soon as a value is written into it, so explicitly reset the alignment
#    back to the worksheet default (General / Bottom / no wrap) to match the
#    plain, un-styled look the other "naked" rows in this table use.
# ---------------------------------------------------------------------------
Set-RowValues $ws 179 @(
    "S21",
    "G03",
    "Portfolio improvement guidelines & triage framework",
    "S21_G03_TD001",
    "Document practical portfolio-stabilisation and profit-framework suggestions (A/B/C buckets, bracket use, risk sizing) in pf_improvement_suggestions.md.",
    "Pure documentation/design task; no code changes beyond the new markdown file.",
    "implemented",
    "Provides a reference playbook for using SigmaTrader tools to manage existing losers and structure new trades systematically."
)
$row179 = $ws.Range("A179:H179")
$row179.WrapText = $false
$row179.VerticalAlignment = -4107  # xlVAlignBottom
$row179.HorizontalAlignment = 1    # xlGeneral

# ---------------------------------------------------------------------------
# 8) Scroll position: topLeftCell moved from C161 to C164 (active cell/
#    selection stays at E169, which we restore explicitly afterwards).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 164
$win.ScrollColumn = 3
$ws.Range("E169").Select()
